# Shift a set of flow-chart shapes on slide 1 to the left (most of them by
# the same amount, ~417443 EMU / 32.87 pt), matching the target OOXML diff.
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are expressed in
# points (1 pt = 12700 EMU) and are backed by single-precision floats, so a
# naive EMU/12700 literal can round to the wrong EMU value on save. The
# literals below were chosen so that, after the float32 round-trip, the
# saved <a:off>/<a:ext> EMU values match the target exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Bent Up Arrow 28 (id 29): x 4910882 -> 4493439
$s.Shapes.Item("Bent Up Arrow 28").Left = 353.8141174316406

# Rectangle 4 (id 5): x 2532118 -> 2114675
$s.Shapes.Item("Rectangle 4").Left = 166.50985717773438

# Rectangle 5 (id 6): x 9359348 -> 8941905
$s.Shapes.Item("Rectangle 5").Left = 704.0870361328125

# Rectangle 6 (id 7): x 9359348 -> 8941905
$s.Shapes.Item("Rectangle 6").Left = 704.0870361328125

# Rectangle 7 (id 8): x 6841752 -> 6424309
$s.Shapes.Item("Rectangle 7").Left = 505.8511047363281

# Right Arrow 8 (id 9): x 8636272 -> 8218829
$s.Shapes.Item("Right Arrow 8").Left = 647.1519165039062

# Rectangle 9 (id 10): x 6080582 -> 5663139
$s.Shapes.Item("Rectangle 9").Left = 445.9164733886719

# Right Arrow 10 (id 11): x 8814437 -> 8396994
$s.Shapes.Item("Right Arrow 10").Left = 661.1806640625

# Rectangle 11 (id 12): x 2243898 -> 1826455
$s.Shapes.Item("Rectangle 11").Left = 143.81536865234375

# Right Arrow 12 (id 13): x 3580908 -> 3163465
$s.Shapes.Item("Right Arrow 12").Left = 249.09173583984375

# Rectangle 14 (id 15): x 6864693 -> 6447250
$s.Shapes.Item("Rectangle 14").Left = 507.6575012207031

# Right Arrow 15 (id 16): x 8659213 -> 8241770
$s.Shapes.Item("Right Arrow 15").Left = 648.9583129882812

# Rectangle 16 (id 17): x 2243898 -> 1826455
$s.Shapes.Item("Rectangle 16").Left = 143.81536865234375

# Rectangle 19 (id 20): x 4204671 -> 3594412, y 3248660 -> 3327143 (moved diagonally)
$s.Shapes.Item("Rectangle 19").Left = 283.02459716796875
$s.Shapes.Item("Rectangle 19").Top = 261.9797668457031

# Right Arrow 20 (id 21): x 3512584 -> 3095141
$s.Shapes.Item("Right Arrow 20").Left = 243.71189880371094

# Rectangle 21 (id 22): x 5913084 -> 5495641
$s.Shapes.Item("Rectangle 21").Left = 432.7276611328125

# Right Arrow 22 (id 23): x 5300190 -> 4882747
$s.Shapes.Item("Right Arrow 22").Left = 384.4682922363281

# Right Arrow 23 (id 24): x 8636271 -> 8218828
$s.Shapes.Item("Right Arrow 23").Left = 647.15185546875

# Right Arrow 24 (id 25): x 5257725 -> 4840282
$s.Shapes.Item("Right Arrow 24").Left = 381.12457275390625

# Rectangle 26 (id 27): x 5052235 -> 4634792
$s.Shapes.Item("Rectangle 26").Left = 364.94427490234375

# Right Arrow 27 (id 28): x 3580908 -> 3163465
$s.Shapes.Item("Right Arrow 27").Left = 249.09173583984375

# Bent Up Arrow 25 (id 26): x 5386907 -> 4969464
$s.Shapes.Item("Bent Up Arrow 25").Left = 391.29638671875

# Rectangle 29 (id 30): x 2119933 -> 1702490
$s.Shapes.Item("Rectangle 29").Left = 134.05433654785156

# Right Arrow 30 (id 31): x 3508488 -> 3091045
$s.Shapes.Item("Right Arrow 30").Left = 243.38937377929688

# Right Arrow 31 (id 32): x 3508488 -> 3091045
$s.Shapes.Item("Right Arrow 31").Left = 243.38937377929688

# Right Arrow 32 (id 33): x 646042 -> 477078, width 1473889 -> 1225409 (left edge moved, right edge moved too)
$s.Shapes.Item("Right Arrow 32").Left = 37.5651969909668
$s.Shapes.Item("Right Arrow 32").Width = 96.48889923095703

# Rectangle 33 (id 34): x 733482 -> 316039
$s.Shapes.Item("Rectangle 33").Left = 24.88496208190918
